# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2410
#   *_new  -> *_FV2504
# Then turn the header+data range into a real Excel Table (adds
# xl/tables/table1.xml + autofilter) and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1, columns A:U) -------------------------
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $oldValue = [string]$cell.Value()
    if ($oldValue.EndsWith("_old")) {
        $cell.Value = $oldValue.Substring(0, $oldValue.Length - 4) + "_FV2410"
    } elseif ($oldValue.EndsWith("_new")) {
        $cell.Value = $oldValue.Substring(0, $oldValue.Length - 4) + "_FV2504"
    }
}

# --- 2. Freeze the header row ----------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert the used range into a native Excel table --------------------
$usedRange = $ws.Range("A1:U67")
$table = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $usedRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$table.Name = "Table1"
$table.TableStyle = ""
